$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking strings stored as TEXT in the source
# workbook (t="inlineStr"). Force text entry (NumberFormat '@') so Excel
# doesn't silently coerce them to numbers/lose significant trailing zeros,
# then clear the formatting delta so no stray number-format style lingers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '244.21'
Set-TextValue $ws.Range("D3") '23.19'
Set-TextValue $ws.Range("D4") '5.403'
Set-TextValue $ws.Range("D5") '0.05986'
Set-TextValue $ws.Range("D6") '3.461'
Set-TextValue $ws.Range("D7") '6.507'
Set-TextValue $ws.Range("D8") '0.8164'
Set-TextValue $ws.Range("D9") '0.9208'
Set-TextValue $ws.Range("D10") '0.1411'
Set-TextValue $ws.Range("D11") '0.07380'
Set-TextValue $ws.Range("D12") '0.03230'
Set-TextValue $ws.Range("D13") '0.03050'
Set-TextValue $ws.Range("D14") '0.09358'
Set-TextValue $ws.Range("D15") '3.856'
Set-TextValue $ws.Range("D16") '0.001563'
Set-TextValue $ws.Range("D17") '0.04665'
Set-TextValue $ws.Range("D18") '0.0005937'
Set-TextValue $ws.Range("D19") '0.006077'
Set-TextValue $ws.Range("D20") '0.005011'
Set-TextValue $ws.Range("D21") '0.0009867'
Set-TextValue $ws.Range("D22") '0.00007895'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D23") '3.625'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range("D24") '2.128'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range("D25") '0.3204'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range("D26") '0.1302'
$ws.Range("E26").Value = '25ProBitTokenPROB'
$ws.Range("B27").Value = 'UpBots'
$ws.Range("C27").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue $ws.Range("D27") '0.0002898'
$ws.Range("E27").Value = '26UpBotsUBXT'
Set-TextValue $ws.Range("D40") '0.03927'
Set-TextValue $ws.Range("D41") '0.006253'
Set-TextValue $ws.Range("D42") '0.1076'
Set-TextValue $ws.Range("D43") '0.002618'
Set-TextValue $ws.Range("D44") '0.007113'
Set-TextValue $ws.Range("D45") '0.00005232'
Set-TextValue $ws.Range("D48") '0.9095'
Set-TextValue $ws.Range("D49") '0.002298'
Set-TextValue $ws.Range("D50") '0.00002099'
Set-TextValue $ws.Range("D51") '0.0001999'
